# Insert two blank rows at the top of the sheet, pushing all existing
# content down by two rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$shp = $ws.Shapes.Item(1)

$ws.Rows("1:2").Insert()

# The emulated Shapes collection does not auto-move the picture's
# cell-anchor when rows are inserted above it (real Excel does this via
# Placement = xlMoveAndSize), so shift it down by the height of the two
# new rows (15pt default row height each) ourselves.
$shp.Top = $shp.Top + 30

# Update the view: show zoomed-out (70%) and select the merged
# "Fluent informative copy..." cell block, matching the target state.
$ws.Application.ActiveWindow.Zoom = 70
$ws.Range("D19:D23").Select()
